# Recurring Deposit 16 Test cases
# The "Output" sheet had an extra "verifyactivationdate" / "Not activated"
# verification row that is no longer needed (client is no longer auto-activated
# on creation) - remove that row, which shifts the remaining verification
# rows up by one and lets Excel drop the now-unused shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# Row 2 held: A2 = "verifyactivationdate", B2 = "Not activated"
$ws.Rows.Item(2).Delete()

# Column A on the Output sheet now needs to fit the longest remaining label
# ("verifyclientclassification"); widen it to match.
$ws.Columns.Item(1).ColumnWidth = 19.59

# Land the selection on A10 and make "Output" the active/visible tab
# (previously "Input" was the selected tab).
$ws.Range("A10").Select() | Out-Null
$ws.Activate() | Out-Null
